# Applies the "Jurisdiction" CodeSystem-metadata-property update to the
# Metadata sheet of the Krebsstadium CodeSystem workbook:
#   - refresh the "Date" property value
#   - insert a new "Jurisdiction" property row (with an empty value) right
#     after "Contact", pushing the rows below it down by one
#
# The "Concepts" sheet (sheet 2) is left untouched - it does not change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Bump the "Date" metadata property (row 8, column B) to the new value.
$ws.Range("B8").Value = "2024-09-17T19:55:11+00:00"

# 2. Make room for the new "Jurisdiction" row at row 11 (just below "Contact").
#    Grab the formatting of the current last row (21) and stamp it onto the
#    new last row (22) first, so that when values are shifted down one row
#    the newly-exposed row 22 ends up with the same look as the rest of the
#    table instead of Excel's bare default style.
$ws.Range("A21:B21").Copy()
$ws.Range("A22:B22").PasteSpecial(-4122)   # xlPasteFormats

# 3. Shift the existing "Description" ... "Count" rows (11-21) down to
#    (12-22), working from the bottom up so we never overwrite a row before
#    reading it.
for ($r = 21; $r -ge 11; $r--) {
    $srcA = "A" + $r
    $srcB = "B" + $r
    $dstA = "A" + ($r + 1)
    $dstB = "B" + ($r + 1)
    $ws.Range($dstA).Value = $ws.Range($srcA).Value2
    $ws.Range($dstB).Value = $ws.Range($srcB).Value2
}

# 4. Populate the freshly-vacated row 11 with the new "Jurisdiction" property
#    (its value is empty, matching the source CodeSystem having no
#    jurisdiction codes).
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
